$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = " formé de béton et d’armatures en acie"
$ws.Range("F2").Value = "σ = N / A"
$ws.Range("G2").Value = "σmax = M × y / I"
$ws.Range("E2").Value = "M = σ × W (ou M = E × I × κ). Ia"

$ws.Range("F5").Select()
